$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.8252008557319641
$ws.Range("B1").Value = 1.978493690490723
$ws.Range("C1").Value = 15
$ws.Range("D1").Value = 1.780203938484192
$ws.Range("E1").Value = 0.4762255549430847
